$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly handled by writing to the new range A2:T11

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd34"
$ws.Range("C2").Value = "Sele"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 181.777022
$ws.Range("H2").Value = 545.331066
$ws.Range("I2").Value = 0.674524008100009
$ws.Range("J2").Value = 0.6869174838889931
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.112632333333332
$ws.Range("N2").Value = 27.337897
$ws.Range("O2").Value = 0.9981738658344552
$ws.Range("P2").Value = 0.9981738658344552
$ws.Range("Q2").Value = 1656.467168134244
$ws.Range("R2").Value = 14908.2045132082
$ws.Range("S2").Value = 0.6732922367633374
$ws.Range("T2").Value = 0.6856630804027534

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cd34"
$ws.Range("C3").Value = "Sele"
$ws.Range("D3").Value = "M2"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 181.777022
$ws.Range("H3").Value = 545.331066
$ws.Range("I3").Value = 0.674524008100009
$ws.Range("J3").Value = 0.6869174838889931
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01667133333333333
$ws.Range("N3").Value = 0.050014
$ws.Range("O3").Value = 0.001826134165544791
$ws.Range("P3").Value = 0.001826134165544791
$ws.Range("Q3").Value = 3.030465326102667
$ws.Range("R3").Value = 27.274187934924
$ws.Range("S3").Value = 0.001231771336671638
$ws.Range("T3").Value = 0.001254403486239754

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cd34"
$ws.Range("C4").Value = "Sele"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 72.32699966666667
$ws.Range("H4").Value = 216.980999
$ws.Range("I4").Value = 0.2683853942167015
$ws.Range("J4").Value = 0.2733166165978157
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 9.112632333333332
$ws.Range("N4").Value = 27.337897
$ws.Range("O4").Value = 0.9981738658344552
$ws.Range("P4").Value = 0.9981738658344552
$ws.Range("Q4").Value = 659.0893557354558
$ws.Range("R4").Value = 5931.804201619102
$ws.Range("S4").Value = 0.2678952864787892
$ws.Range("T4").Value = 0.2728175037862353

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cd34"
$ws.Range("C5").Value = "Sele"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 72.32699966666667
$ws.Range("H5").Value = 216.980999
$ws.Range("I5").Value = 0.2683853942167015
$ws.Range("J5").Value = 0.2733166165978157
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01667133333333333
$ws.Range("N5").Value = 0.050014
$ws.Range("O5").Value = 0.001826134165544791
$ws.Range("P5").Value = 0.001826134165544791
$ws.Range("Q5").Value = 1.205787520442889
$ws.Range("R5").Value = 10.852087683986
$ws.Range("S5").Value = 0.000490107737912326
$ws.Range("T5").Value = 0.0004991128115803778

# Row 6
$ws.Range("A6").Value = "M1"
$ws.Range("B6").Value = "Cd34"
$ws.Range("C6").Value = "Sele"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.1272716666666667
$ws.Range("H6").Value = 0.381815
$ws.Range("I6").Value = 0.0004722697829078107
$ws.Range("J6").Value = 0.0004809471080290077
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.112632333333332
$ws.Range("N6").Value = 27.337897
$ws.Range("O6").Value = 0.9981738658344552
$ws.Range("P6").Value = 0.9981738658344552
$ws.Range("Q6").Value = 1.159779904783889
$ws.Range("R6").Value = 10.438019143055
$ws.Range("S6").Value = 0.0004714073549218883
$ws.Range("T6").Value = 0.000480068834083216

# Row 7
$ws.Range("A7").Value = "M1"
$ws.Range("B7").Value = "Cd34"
$ws.Range("C7").Value = "Sele"
$ws.Range("D7").Value = "M2"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.1272716666666667
$ws.Range("H7").Value = 0.381815
$ws.Range("I7").Value = 0.0004722697829078107
$ws.Range("J7").Value = 0.0004809471080290077
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01667133333333333
$ws.Range("N7").Value = 0.050014
$ws.Range("O7").Value = 0.001826134165544791
$ws.Range("P7").Value = 0.001826134165544791
$ws.Range("Q7").Value = 0.002121788378888889
$ws.Range("R7").Value = 0.01909609541
$ws.Range("S7").Value = [double]"8.624279859223746e-07"
$ws.Range("T7").Value = [double]"8.782739457917324e-07"

# Row 8
$ws.Range("A8").Value = "M2"
$ws.Range("B8").Value = "Cd34"
$ws.Range("C8").Value = "Sele"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.6715256666666667
$ws.Range("H8").Value = 2.014577
$ws.Range("I8").Value = 0.002491845114626373
$ws.Range("J8").Value = 0.002537629433237966
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.112632333333332
$ws.Range("N8").Value = 27.337897
$ws.Range("O8").Value = 0.9981738658344552
$ws.Range("P8").Value = 0.9981738658344552
$ws.Range("Q8").Value = 6.119366502729888
$ws.Range("R8").Value = 55.074298524569
$ws.Range("S8").Value = 0.002487294671127308
$ws.Range("T8").Value = 0.002532995381430439

# Row 9
$ws.Range("A9").Value = "M2"
$ws.Range("B9").Value = "Cd34"
$ws.Range("C9").Value = "Sele"
$ws.Range("D9").Value = "M2"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.6715256666666667
$ws.Range("H9").Value = 2.014577
$ws.Range("I9").Value = 0.002491845114626373
$ws.Range("J9").Value = 0.002537629433237966
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01667133333333333
$ws.Range("N9").Value = 0.050014
$ws.Range("O9").Value = 0.001826134165544791
$ws.Range("P9").Value = 0.001826134165544791
$ws.Range("Q9").Value = 0.01119522823088889
$ws.Range("R9").Value = 0.100757054078
$ws.Range("S9").Value = [double]"4.550443499065095e-06"
$ws.Range("T9").Value = [double]"4.634051807527915e-06"

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cd34"
$ws.Range("C10").Value = "Sele"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.5865095
$ws.Range("H10").Value = 29.173019
$ws.Range("I10").Value = 0.05412648278575528
$ws.Range("J10").Value = 0.03674732297192435
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.112632333333332
$ws.Range("N10").Value = 27.337897
$ws.Range("O10").Value = 0.9981738658344552
$ws.Range("P10").Value = 0.9981738658344552
$ws.Range("Q10").Value = 132.9214981001738
$ws.Range("R10").Value = 797.5289886010429
$ws.Range("S10").Value = 0.05402764056627944
$ws.Range("T10").Value = 0.036680217429953

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Cd34"
$ws.Range("C11").Value = "Sele"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 14.5865095
$ws.Range("H11").Value = 29.173019
$ws.Range("I11").Value = 0.05412648278575528
$ws.Range("J11").Value = 0.03674732297192435
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01667133333333333
$ws.Range("N11").Value = 0.050014
$ws.Range("O11").Value = 0.001826134165544791
$ws.Range("P11").Value = 0.001826134165544791
$ws.Range("Q11").Value = 0.2431765620443333
$ws.Range("R11").Value = 1.459059372266
$ws.Range("S11").Value = [double]"9.884221947583974e-05"
$ws.Range("T11").Value = [double]"6.710554197134e-05"

